$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3566.6667
$ws.Range("I6").Value = 350
$ws.Range("K6").Value = 1050
$ws.Range("M6").Value = -938

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1083.3334
$ws.Range("I8").Value = 125
$ws.Range("K8").Value = 375
$ws.Range("M8").Value = -236

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 100031.45
$ws.Range("I11").Value = 100031.45
$ws.Range("K11").Value = 100031.45
$ws.Range("M11").Value = -99891.45

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2400
$ws.Range("I40").Value = 1900
$ws.Range("J40").Value = 2700
$ws.Range("K40").Value = 1900
$ws.Range("L40").Value = 2700
$ws.Range("M40").Value = -1725
$ws.Range("N40").Value = -3050

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 740.59375
$ws.Range("I98").Value = 670.44446
$ws.Range("K98").Value = 670.44446
$ws.Range("M98").Value = 827.55554

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 76927710
$ws.Range("I113").Value = 333335840
$ws.Range("J113").Value = 5280.3
$ws.Range("K113").Value = 333335840
$ws.Range("L113").Value = 5280.3
$ws.Range("M113").Value = -333332586
$ws.Range("N113").Value = -11788.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 740.59375
$ws.Range("I122").Value = 670.44446
$ws.Range("K122").Value = 2011.33338
$ws.Range("M122").Value = 438.66662

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1879.9231
$ws.Range("I138").Value = 1418.931
$ws.Range("J138").Value = 3216.8
$ws.Range("K138").Value = 4256.793
$ws.Range("L138").Value = 9650.400000000001
$ws.Range("M138").Value = 883.2070000000003
$ws.Range("N138").Value = -19930.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2876.5208
$ws.Range("I32").Value = 2405.8086
$ws.Range("K32").Value = 2405.8086
$ws.Range("M32").Value = -2118.8086

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 43000
$ws.Range("J56").Value = 43000
$ws.Range("L56").Value = 43000
$ws.Range("N56").Value = -44484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3191.875
$ws.Range("I61").Value = 3076.4285
$ws.Range("K61").Value = 3076.4285
$ws.Range("M61").Value = -2864.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3765
$ws.Range("I122").Value = 2941.6667
$ws.Range("K122").Value = 8825.000100000001
$ws.Range("M122").Value = -6375.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3191.875
$ws.Range("I136").Value = 3076.4285
$ws.Range("K136").Value = 9229.2855
$ws.Range("M136").Value = -6679.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 69999.336
$ws.Range("J38").Value = 69999.336
$ws.Range("L38").Value = 69999.336
$ws.Range("N38").Value = -70831.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 16674333
$ws.Range("I44").Value = 25001500
$ws.Range("J44").Value = 20000
$ws.Range("K44").Value = 25001500
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = -25001003
$ws.Range("N44").Value = -20994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1819
$ws.Range("I105").Value = 1621.1875
$ws.Range("K105").Value = 1621.1875
$ws.Range("M105").Value = 125.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 53882.668
$ws.Range("J50").Value = 65824.664
$ws.Range("L50").Value = 65824.664
$ws.Range("N50").Value = -67074.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2750
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -1002
$ws.Range("N99").Value = -6496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1800.6957
$ws.Range("I134").Value = 1405.6316
$ws.Range("J134").Value = 3677.25
$ws.Range("K134").Value = 4216.8948
$ws.Range("L134").Value = 11031.75
$ws.Range("M134").Value = -1681.8948
$ws.Range("N134").Value = -16101.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 4999
$ws.Range("I43").Value = 4999
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 14997
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -14883
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1021.625
$ws.Range("J97").Value = 1164.6
$ws.Range("L97").Value = 3493.8
$ws.Range("N97").Value = -4485.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1510.45
$ws.Range("J113").Value = 1675.5625
$ws.Range("L113").Value = 5026.6875
$ws.Range("N113").Value = -9366.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1973.5454
$ws.Range("I137").Value = 1164
$ws.Range("J137").Value = 2648.1667
$ws.Range("K137").Value = 3492
$ws.Range("L137").Value = 7944.500100000001
$ws.Range("M137").Value = 1608
$ws.Range("N137").Value = -18144.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 21749.8
$ws.Range("J98").Value = 21749.8
$ws.Range("L98").Value = 21749.8
$ws.Range("N98").Value = -27739.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1552.0731
$ws.Range("I102").Value = 1341.5135
$ws.Range("J102").Value = 3499.75
$ws.Range("K102").Value = 1341.5135
$ws.Range("L102").Value = 3499.75
$ws.Range("M102").Value = 280.4865
$ws.Range("N102").Value = -6743.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4613544
$ws.Range("I22").Value = 345.66666
$ws.Range("K22").Value = 345.66666
$ws.Range("M22").Value = -50.66665999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4613544
$ws.Range("I27").Value = 345.66666
$ws.Range("K27").Value = 345.66666
$ws.Range("M27").Value = -238.66666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1826.375
$ws.Range("I46").Value = 1597.6875
$ws.Range("J46").Value = 2283.75
$ws.Range("K46").Value = 1597.6875
$ws.Range("L46").Value = 2283.75
$ws.Range("M46").Value = -1409.6875
$ws.Range("N46").Value = -2659.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 37354.332
$ws.Range("J54").Value = 35312.715
$ws.Range("L54").Value = 35312.715
$ws.Range("N54").Value = -36600.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2464.4119
$ws.Range("I136").Value = 1955.7778
$ws.Range("J136").Value = 4426.2856
$ws.Range("K136").Value = 5867.3334
$ws.Range("L136").Value = 13278.8568
$ws.Range("M136").Value = -3317.3334
$ws.Range("N136").Value = -18378.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 4699.857
$ws.Range("J18").Value = 13949.5
$ws.Range("L18").Value = 13949.5
$ws.Range("N18").Value = -14295.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14297287
$ws.Range("I81").Value = 10399.2
$ws.Range("J81").Value = 22234446
$ws.Range("K81").Value = 20798.4
$ws.Range("L81").Value = 44468892
$ws.Range("M81").Value = -19737.4
$ws.Range("N81").Value = -44471014

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 14297287
$ws.Range("I84").Value = 10399.2
$ws.Range("J84").Value = 22234446
$ws.Range("K84").Value = 103992
$ws.Range("L84").Value = 222344460
$ws.Range("M84").Value = -98688
$ws.Range("N84").Value = -222355068

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3088.375
$ws.Range("I100").Value = 3491.8
$ws.Range("J100").Value = 1071.25
$ws.Range("K100").Value = 6983.6
$ws.Range("L100").Value = 2142.5
$ws.Range("M100").Value = -6442.6
$ws.Range("N100").Value = -3224.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2702.9119
$ws.Range("I136").Value = 1015.7308
$ws.Range("J136").Value = 8186.25
$ws.Range("K136").Value = 3047.1924
$ws.Range("L136").Value = 24558.75
$ws.Range("M136").Value = -497.1923999999999
$ws.Range("N136").Value = -29658.75
